$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.301.30"
$ws.Cells.Item(2, 5).Value = "  -2.78%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.571.24"
$ws.Cells.Item(3, 5).Value = "  -3.81%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.20%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'207.54"
$ws.Cells.Item(5, 5).Value = "  -3.26%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.17%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.477"
$ws.Cells.Item(7, 5).Value = "  -4.93%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -2.33%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.0607"
$ws.Cells.Item(9, 5).Value = "  -1.93%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -2.20%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.83%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.795.27"
$ws.Cells.Item(12, 5).Value = "  -3.46%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.580.64"
$ws.Cells.Item(13, 5).Value = "  -3.16%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -3.54%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.505"
$ws.Cells.Item(15, 5).Value = "  -3.85%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "25.309.63"
$ws.Cells.Item(16, 5).Value = "  -2.66%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'59.71"
$ws.Cells.Item(17, 5).Value = "  -2.73%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "0.0₃0707"
$ws.Cells.Item(18, 5).Value = "  -4.48%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.03%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'185.67"
$ws.Cells.Item(20, 5).Value = "  -2.68%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -2.40%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -3.35%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'5.89"
$ws.Cells.Item(23, 5).Value = "  -2.98%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.10%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'141.16"
$ws.Cells.Item(25, 5).Value = "  -2.15%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -3.04%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -4.81%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'14.86"
$ws.Cells.Item(28, 5).Value = "  -2.10%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'6.44"
$ws.Cells.Item(29, 5).Value = "  -4.78%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -5.65%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.0460"
$ws.Cells.Item(31, 5).Value = "  -4.22%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.05"
$ws.Cells.Item(32, 5).Value = "  -2.46%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'3.00"
$ws.Cells.Item(33, 5).Value = "  -4.01%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "LidoDAOToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(34, 4).Value = "'1.46"
$ws.Cells.Item(34, 5).Value = "  -1.57%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).Value = "'2.26"
$ws.Cells.Item(35, 5).Value = "  -6.40%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "1.084.61"
$ws.Cells.Item(36, 5).Value = "  -3.94%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'1.01"
$ws.Cells.Item(37, 5).Value = "  -0.36%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'2.33"
$ws.Cells.Item(38, 5).Value = "  -4.40%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -2.46%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.778"
$ws.Cells.Item(40, 5).Value = "  -9.59%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.492"
$ws.Cells.Item(41, 5).Value = "  -4.88%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "'0.760"
$ws.Cells.Item(42, 5).Value = "  -1.95%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).Value = "'93.49"
$ws.Cells.Item(43, 5).Value = "  -4.91%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value = "'5.05"
$ws.Cells.Item(44, 5).Value = "  -3.25%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "RocketPoolETH"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(45, 4).Value = "1.707.36"
$ws.Cells.Item(45, 5).Value = "  -3.53%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.0₆0106"
$ws.Cells.Item(46, 5).Value = "  -8.13%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'52.67"
$ws.Cells.Item(47, 5).Value = "  -3.77%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0508"
$ws.Cells.Item(48, 5).Value = "  -3.65%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -1.61%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "USDD"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(50, 4).Value = "'1.01"
$ws.Cells.Item(50, 5).Value = "  -0.12%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).Value = "'1.40"
$ws.Cells.Item(51, 5).Value = "  -5.58%  "
